# Update tags and description
# - The six "Places" entries p140..p145 (column G, rows 82-87) were removed
#   (folded into "PLC2" / no longer listed individually).
# - Three new Transition/Boundary pairs were appended: t150/b150, t151/b151,
#   t152/b152 (columns H/I, new rows 91-93).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete Places rows (clears the cells entirely, shrinking
# the used range the same way Excel does when you select & delete content).
$ws.Range("G82:G87").ClearContents()

# Append the three new Transition / Boundary entries.
$ws.Range("H91").Value = "t150"
$ws.Range("I91").Value = "b150"
$ws.Range("H92").Value = "t151"
$ws.Range("I92").Value = "b151"
$ws.Range("H93").Value = "t152"
$ws.Range("I93").Value = "b152"

# Leave the view roughly where the author left it (cosmetic only).
$ws.Range("G82:G87").Select()
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 78
    $win.ScrollColumn = 1
}
